$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column C ("attributes"), shifting attributes -> D
# and sort_order -> E. This gives room for the new "label-fr_FR" column.
$ws.Columns("C:C").Insert()

# New header for the inserted column
$ws.Range("C1").Value = "label-fr_FR"

# French label for the "other" attribute group
$ws.Range("C2").Value = "Autre"

# erp: add PACK_SIZE to the attributes list
$ws.Range("D3").Value = $ws.Range("D3").Value2 + ",PACK_SIZE"

# PRODUCT_CODING: add WHO's COVID-19 product list related attributes
$ws.Range("D5").Value = $ws.Range("D5").Value2 + ",LMIS_CODE,LMIS_UUID,WHO_COVID19_ITEM_CODE,UN_SUPPLY_MATERIAL_CODE"

# New row 31: _LMIS attribute group
$ws.Range("A31").Value = "_LMIS"
$ws.Range("B31").Value = "LMIS"
$ws.Range("C31").Value = "FR LMIS"
$ws.Range("D31").Value = "PACK_ROUNDING_THRESHOLD,LMIS_ROUND_TO_ZERO"
$ws.Range("E31").Value = "'29"
